$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 104778214
$ws.Range("B3").Value = 77506
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 346435.7253240386
$ws.Range("R3").Value = 6587685.342734368
$ws.Range("AC3").Value = "På tallstam"

# Row 4
$ws.Range("A4").Value = 104778226
$ws.Range("B4").Value = 90653
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 4364
$ws.Range("F4").Value = "Dropptaggsvamp"
$ws.Range("G4").Value = "Hydnellum ferrugineum"
$ws.Range("H4").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q4").Value = 346347.7517754274
$ws.Range("R4").Value = 6587694.528751616
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("A5").Value = 104778356
$ws.Range("B5").Value = 94121
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = "Vedtrappmossa"
$ws.Range("G5").Value = "Crossocalyx hellerianus"
$ws.Range("H5").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q5").Value = 346597.3486557564
$ws.Range("R5").Value = 6587540.61948923
$ws.Range("AC5").Value = "Granlåga"

# Row 6
$ws.Range("A6").Value = 104778248
$ws.Range("B6").Value = 94121
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 53
$ws.Range("F6").Value = "Vedtrappmossa"
$ws.Range("G6").Value = "Crossocalyx hellerianus"
$ws.Range("H6").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q6").Value = 346635.7879304852
$ws.Range("R6").Value = 6587418.232812922
$ws.Range("AC6").Value = "På granlåga"

# Row 7
$ws.Range("A7").Value = 104778340
$ws.Range("B7").Value = 73631
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6426
$ws.Range("F7").Value = "Kattfotslav"
$ws.Range("G7").Value = "Felipes leucopellaeus"
$ws.Range("H7").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q7").Value = 346631.1643096056
$ws.Range("R7").Value = 6587581.049807825
$ws.Range("AC7").ClearContents()

# Row 8
$ws.Range("A8").Value = 104778308
$ws.Range("B8").Value = 89356
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 5447
$ws.Range("F8").Value = "Vedticka"
$ws.Range("G8").Value = "Fuscoporia viticola"
$ws.Range("H8").Value = "(Schwein.) Murrill"
$ws.Range("Q8").Value = 346760.5213297271
$ws.Range("R8").Value = 6587773.595144214
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 104778261
$ws.Range("B9").Value = 77506
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 346708.7430858334
$ws.Range("R9").Value = 6587617.664139647
$ws.Range("AC9").Value = "Rikligt på flera tallstammar i gammal hällmarksskog"

# Row 10
$ws.Range("A10").Value = 104778295
$ws.Range("B10").Value = 94121
$ws.Range("E10").Value = 53
$ws.Range("F10").Value = "Vedtrappmossa"
$ws.Range("G10").Value = "Crossocalyx hellerianus"
$ws.Range("H10").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q10").Value = 346760.5213297271
$ws.Range("R10").Value = 6587773.595144214
$ws.Range("AC10").Value = "Granlåga i källdråg"

# Row 11
$ws.Range("A11").Value = 104778334
$ws.Range("B11").Value = 77506
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 346631.1643096056
$ws.Range("R11").Value = 6587581.049807825
$ws.Range("AC11").Value = "På gran"
